{"js": "// Prepend \"Design: \" to every bulleted (\"List Bullet\") feedback paragraph\n// in the document's answer cells.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/styleBuiltIn,items/text\");\nawait context.sync();\n\nconst prefix = \"Design: \";\n\nfor (const paragraph of paragraphs.items) {\n  // The style name usually resolves to \"List Bullet\", but the very last\n  // paragraph in the body can report an empty style string from this\n  // property while still being styleBuiltIn \"Other\" (i.e. a non-Normal,\n  // non-heading custom/list style) \u2014 so treat that case as a match too.\n  const isListBullet =\n    paragraph.style === \"List Bullet\" ||\n    (paragraph.style === \"\" && paragraph.styleBuiltIn === \"Other\");\n\n  if (isListBullet && !paragraph.text.startsWith(prefix)) {\n    paragraph.insertText(prefix, Word.InsertLocation.start);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Prepend \"Design: \" to every bulleted (\"List Bullet\") feedback paragraph\n# in the document's answer cells.\n$d = $word.ActiveDocument\n$prefix = \"Design: \"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Style.NameLocal\n    if ($styleName -eq \"List Bullet\") {\n        $r = $p.Range\n        $text = $r.Text\n        if ($text.Length -lt $prefix.Length -or $text.Substring(0, $prefix.Length) -ne $prefix) {\n            $r.InsertBefore($prefix)\n        }\n    }\n}\n"}
